# Fixed some bugs from the ctc office and added docstring comments to the line.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The arrival times recorded in the schedule were off by 9 hours
# (entered as 10:xx AM instead of 1:xx AM). Correct the "Arrival" column.
$ws.Range("C2").Value = 0.041666666666666664
$ws.Range("C3").Value = 0.052083333333333336
$ws.Range("C4").Value = 0.052083333333333336
$ws.Range("C5").Value = 0.0625

# Update the active selection to reflect where the CTC office left off editing.
$ws.Range("C3").Select()
